$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the match-data columns (F:V) between pairs of rows. ---
# Columns A:E (Indice, pais, torneio, temporada, data_partida) stay put;
# only the match info (home/away teams, goals, odds, timestamps, url)
# moves, which is what the diff shows for each pair below.
$swapPairs = @(
    @(2, 3),
    @(6, 7),
    @(18, 19),
    @(28, 29),
    @(30, 31),
    @(48, 49),
    @(53, 55),
    @(59, 60),
    @(78, 79),
    @(86, 87),
    @(126, 127)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("F" + $r1 + ":V" + $r1)
    $range2 = $ws.Range("F" + $r2 + ":V" + $r2)
    $vals1 = $range1.Value()
    $vals2 = $range2.Value()
    $range1.Value = $vals2
    $range2.Value = $vals1
}

# --- 2) Append the two new match rows at the bottom (134, 135). ---
$newRows = @(
    @(133, "italy", "serie-a", "2023-2024", 45262.75, "Lazio", 1, "Cagliari", 0,
      1.58, "15/11/2023 16:01", 1.74, "02/12/2023 17:58",
      4.08, "15/11/2023 16:01", 3.84, "02/12/2023 17:58",
      5.96, "15/11/2023 16:01", 5.16, "02/12/2023 17:58",
      "https://www.betexplorer.com/football/italy/serie-a/lazio-cagliari/tt1JVBW2/"),
    @(134, "italy", "serie-a", "2023-2024", 45262.86458333334, "AC Milan", 3, "Frosinone", 1,
      1.31, "15/11/2023 16:01", 1.45, "02/12/2023 20:39",
      5.9, "15/11/2023 16:01", 4.88, "02/12/2023 20:44",
      10.28, "15/11/2023 16:01", 7.17, "02/12/2023 20:44",
      "https://www.betexplorer.com/football/italy/serie-a/ac-milan-frosinone/x8fSTk1F/")
)

$targetRow = 134
foreach ($rowData in $newRows) {
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $rowData[$c]
    }
    $targetRow = $targetRow + 1
}

# Copy the Indice (A) / data_partida (E) cell formatting from the last
# pre-existing data row (133) so the new rows match the sheet's styling
# (bold+border+centered index column, date-time number format column).
$ws.Range("A133").Copy()
$ws.Range("A134:A135").PasteSpecial(-4122)

$ws.Range("E133").Copy()
$ws.Range("E134:E135").PasteSpecial(-4122)
